$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1294
$ws.Range("F5").Value = 104
$ws.Range("F6").Value = 365
$ws.Range("F7").Value = 1188
$ws.Range("F8").Value = 446
$ws.Range("F9").Value = 7139
$ws.Range("F11").Value = 91
$ws.Range("F12").Value = 2049
$ws.Range("F13").Value = 8003
$ws.Range("F16").Value = 5517
$ws.Range("F18").Value = 2431
$ws.Range("F19").Value = 1030
$ws.Range("F21").Value = 305
$ws.Range("F23").Value = 79
$ws.Range("F25").Value = 385
$ws.Range("F28").Value = 2388
$ws.Range("F30").Value = 266
$ws.Range("F31").Value = 85
$ws.Range("F32").Value = 156
$ws.Range("F33").Value = 588
$ws.Range("F36").Value = 1514
$ws.Range("F39").Value = 2346
$ws.Range("F40").Value = 2219
$ws.Range("F41").Value = 2
$ws.Range("F42").Value = 10

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 81
$ws.Range("F4").Value = 69
$ws.Range("F5").Value = 9
$ws.Range("F6").Value = 28

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 257

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 257
$ws.Range("F6").Value = 1294
$ws.Range("F8").Value = 365
$ws.Range("F9").Value = 1188
$ws.Range("F10").Value = 446
$ws.Range("F11").Value = 7139
$ws.Range("F13").Value = 91
$ws.Range("F14").Value = 2049
$ws.Range("F15").Value = 8003
$ws.Range("F18").Value = 5517
$ws.Range("F20").Value = 2431
$ws.Range("F21").Value = 1030
$ws.Range("F24").Value = 79
$ws.Range("F25").Value = 81
$ws.Range("F27").Value = 69
$ws.Range("F28").Value = 385
$ws.Range("F30").Value = 2388
$ws.Range("F32").Value = 266
$ws.Range("F33").Value = 85
$ws.Range("F34").Value = 156
$ws.Range("F35").Value = 9
$ws.Range("F36").Value = 588
$ws.Range("F39").Value = 28
$ws.Range("F40").Value = 1514
$ws.Range("F43").Value = 2346
$ws.Range("F45").Value = 2219
$ws.Range("F46").Value = 2
$ws.Range("F47").Value = 10
